$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 38: new activity entry
$ws.Range("A38").Value = 43218
$ws.Range("B38").Value = "Transfère des données du compte bancaire du controlleur bankAccount au controlleur detailBankAccount. Quelque modification dans la fenêtre du détail d'un compte bancaire.  Modification des champs et test des graphiques"
$ws.Range("C38").Value = 2

# Row 39: new activity entry
$ws.Range("A39").Value = 43218
$ws.Range("B39").Value = "Rapport sur la partie détail compte bancaire"
$ws.Range("C39").Value = 0.5

# Update the view state: scroll position and selection to mirror Excel's UI state after edit
$ws.Range("G39").Select()
